$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new appropriation entry (row 13): date 2013-10-18, 1h30min worked
# (0.0625 of a day = 1.5 hours). Reuse the date/time number formats already
# used by the rows above by copying their cell formatting onto the new cells
# before writing the values, so the new row renders like the existing ones.
$ws.Range("A12").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value = 41565

$ws.Range("B12").Copy()
$ws.Range("B13").PasteSpecial(-4122)
$ws.Range("B13").Value = 0.0625

# Match the author's new selection after entering the data.
$ws.Range("C13").Select()
